$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 94, shifting existing rows 94-174 down to 95-175
$ws.Rows.Item(94).Insert()

$ws.Cells.Item(94, 1).Value = 4
$ws.Cells.Item(94, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(94, 3).Value = "Los Lagos"
$ws.Cells.Item(94, 4).Value = "2023-08-29"
$ws.Cells.Item(94, 5).Value = 10
$ws.Cells.Item(94, 6).Value = 100112022
$ws.Cells.Item(94, 7).Value = "Arveja Verde"
$ws.Cells.Item(94, 8).Value = "Perfection"
$ws.Cells.Item(94, 9).Value = "Primera"
$ws.Cells.Item(94, 10).Value = 50
$ws.Cells.Item(94, 11).Value = 38000
$ws.Cells.Item(94, 12).Value = 38000
$ws.Cells.Item(94, 13).Value = 38000
$ws.Cells.Item(94, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(94, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(94, 16).Value = 1520
$ws.Cells.Item(94, 17).Value = 25
$ws.Cells.Item(94, 18).Value = "Hortaliza"
